$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.225.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.57%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.010.18'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').Value = '  +1.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.73'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.70%  '
$ws.Range('E9').Value = '  +3.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0811'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.61%  '
$ws.Range('E13').Value = '  +2.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.304.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.38%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('E16').Value = '  +3.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.019.40'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.115.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0866'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('E21').Value = '  +2.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.83'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.27%  '
$ws.Range('E28').Value = '  -2.52%  '
$ws.Range('E29').Value = '  +1.33%  '
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.82'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('E33').Value = '  +5.86%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.47'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.48'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.70%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E38').Value = '  +1.92%  '
$ws.Range('E39').Value = '  -4.61%  '
$ws.Range('E40').Value = '  -0.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.92'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('E43').Value = '  +1.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.65'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.374.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('E47').Value = '  +1.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.43'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.61%  '
$ws.Range('E49').Value = '  +12.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.85'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.85%  '
